# Bayes-1.xlsx — rebuild the sheet with the new "Pocasi / Jdeme behat?" table
# (2 columns x 5 rows) in place of the old 5-column weather table, rename the
# sheet from "List1" to "Bayes-1", and leave the selection on B6.

$wb = $excel.ActiveWorkbook

$oldName = $wb.ActiveSheet.Name

# Build the new sheet from scratch (rather than editing the old one in
# place) so the worksheet doesn't keep stale <col> width metadata for the
# columns (A, C, D, E) that no longer carry any custom formatting.
$new = $wb.Worksheets.Add()
$old = $wb.Worksheets.Item($oldName)

$new.Range("A1").Value = "Počasí"
$new.Range("B1").Value = "Jdeme běhat?"

$new.Range("A2").Value = "slunečno"
$new.Range("B2").Value = "ne"

$new.Range("A3").Value = "slunečno"
$new.Range("B3").Value = "ano"

$new.Range("A4").Value = "zataženo "
$new.Range("B4").Value = "ano"

$new.Range("A5").Value = "déšť"
$new.Range("B5").Value = "ne"

# Widen column B to fit the new header/values.
$new.Columns("B").ColumnWidth = 13.71

# Drop the old sheet and take its place under the new name.
$old.Delete() | Out-Null
$new.Name = "Bayes-1"

$new.Range("B6").Select() | Out-Null
